$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data table (rows 2-11): 5 sending clusters x 2 target clusters (FAPs, MuSCs).
# The previous "ECs" target-cluster rows are gone and a new "Inflammatory-Mac"
# sending cluster was added, per the updated TPM run.
$rows = @(
    @{row=2; vals=@('ECs', 'Ncam1', 'Ptprz1', 'FAPs', 3, 1, 2.721912, 8.165736000000001, 0.03958906885744382, 0.03958906885744382, 2, 0.6666666666666666, 0.05372733333333333, 0.161182, 0.1072370469527173, 0.1072370469527173, 0.146241073328, 1.316169659952, 0.00424541483588006, 0.004245414835880062)},
    @{row=3; vals=@('ECs', 'Ncam1', 'Ptprz1', 'MuSCs', 3, 1, 2.721912, 8.165736000000001, 0.03958906885744382, 0.03958906885744382, 3, 1, 0.4472873333333333, 1.341862, 0.8927629530472828, 0.8927629530472827, 1.217476760048, 10.957290840432, 0.03534365402156376, 0.03534365402156376)},
    @{row=4; vals=@('FAPs', 'Ncam1', 'Ptprz1', 'FAPs', 3, 1, 22.91769833333333, 68.753095, 0.3333283138369127, 0.3333283138369127, 2, 0.6666666666666666, 0.05372733333333333, 0.161182, 0.1072370469527173, 0.1072370469527173, 1.231306817587778, 11.08176135829, 0.03574514404159909, 0.03574514404159909)},
    @{row=5; vals=@('FAPs', 'Ncam1', 'Ptprz1', 'MuSCs', 3, 1, 22.91769833333333, 68.753095, 0.3333283138369127, 0.3333283138369127, 3, 1, 0.4472873333333333, 1.341862, 0.8927629530472828, 0.8927629530472827, 10.25079617365444, 92.25716556288999, 0.2975831697953136, 0.2975831697953136)},
    @{row=6; vals=@('Inflammatory-Mac', 'Ncam1', 'Ptprz1', 'FAPs', 1, 0.3333333333333333, 0.02880066666666667, 0.08640200000000001, 0.0004188936217654919, 0.0004188936217654919, 2, 0.6666666666666666, 0.05372733333333333, 0.161182, 0.1072370469527173, 0.1072370469527173, 0.001547383018222222, 0.013926447164, [double]"4.492091498545985E-05", [double]"4.492091498545985E-05")},
    @{row=7; vals=@('Inflammatory-Mac', 'Ncam1', 'Ptprz1', 'MuSCs', 1, 0.3333333333333333, 0.02880066666666667, 0.08640200000000001, 0.0004188936217654919, 0.0004188936217654919, 3, 1, 0.4472873333333333, 1.341862, 0.8927629530472828, 0.8927629530472827, 0.01288217339155556, 0.115939560524, 0.0003739727067800321, 0.000373972706780032)},
    @{row=8; vals=@('MuSCs', 'Ncam1', 'Ptprz1', 'FAPs', 3, 1, 43.03445633333334, 129.103369, 0.6259181248412854, 0.6259181248412854, 2, 0.6666666666666666, 0.05372733333333333, 0.161182, 0.1072370469527173, 0.1072370469527173, 2.312126580239778, 20.809139222158, 0.06712161134216169, 0.06712161134216169)},
    @{row=9; vals=@('MuSCs', 'Ncam1', 'Ptprz1', 'MuSCs', 3, 1, 43.03445633333334, 129.103369, 0.6259181248412854, 0.6259181248412854, 3, 1, 0.4472873333333333, 1.341862, 0.8927629530472828, 0.8927629530472827, 19.24876721478645, 173.238904933078, 0.5587965134991237, 0.5587965134991237)},
    @{row=10; vals=@('Resolving-Mac', 'Ncam1', 'Ptprz1', 'FAPs', 1, 0.3333333333333333, 0.051263, 0.153789, 0.0007455988425926857, 0.0007455988425926857, 2, 0.6666666666666666, 0.05372733333333333, 0.161182, 0.1072370469527173, 0.1072370469527173, 0.002754224288666667, 0.024788018598, [double]"7.99558180910035E-05", [double]"7.995581809100351E-05")},
    @{row=11; vals=@('Resolving-Mac', 'Ncam1', 'Ptprz1', 'MuSCs', 1, 0.3333333333333333, 0.051263, 0.153789, 0.0007455988425926857, 0.0007455988425926857, 3, 1, 0.4472873333333333, 1.341862, 0.8927629530472828, 0.8927629530472827, 0.02292929056866667, 0.206363615118, 0.0006656430245016823, 0.0006656430245016822)}
)

foreach ($item in $rows) {
    $r = $item.row
    $c = 1
    foreach ($v in $item.vals) {
        $ws.Cells.Item($r, $c).Value = $v
        $c = $c + 1
    }
}

# Drop the now-unused trailing rows (previously 12-13, "Resolving-Mac" / "ECs" pairing
# removed) so the sheet dimension shrinks back down to T11.
$ws.Range("A12:T13").ClearContents()